# Update cryptos list with the latest scraped values
# (mirrors the "Updated cryptos list ... with GitHub Actions" data-refresh commit).
#
# D-column price cells are text (inlineStr) in the source data, e.g. a price
# of "462.00" must keep its trailing zero and NOT turn into the number 462 -
# so NumberFormat is forced to Text ("@") before each D-cell write.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.432.79"
$ws.Range("E2").Value = "  +1.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.173.44"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.74"
$ws.Range("E5").Value = "  +2.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.54"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.172.53"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("E10").Value = "  +0.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.35"
$ws.Range("E11").Value = "  +2.09%  "

$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("E13").Value = "  +2.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.83"
$ws.Range("E14").Value = "  +5.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.695.68"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.169.66"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.400.82"

$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.00"
$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("E22").Value = "  -0.62%  "

$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.28"
$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.18"
$ws.Range("E25").Value = "  +0.91%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("E29").Value = "  +2.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.80"
$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("E31").Value = "  -1.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.18"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("E33").Value = "  -1.42%  "

$ws.Range("E34").Value = "  +1.11%  "

$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("E36").Value = "  +2.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.30"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0731"
$ws.Range("E38").Value = "  +6.50%  "

$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("E40").Value = "  +1.83%  "

$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "393.35"
$ws.Range("E43").Value = "  -4.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.804.32"
$ws.Range("E44").Value = "  -5.04%  "

$ws.Range("E45").Value = "  +0.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.87"
$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.12"
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.61"
$ws.Range("E49").Value = "  +2.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.17"
$ws.Range("E50").Value = "  -1.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.112"
$ws.Range("E51").Value = "  +0.54%  "
